$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trend_instructions")

# New header cells (row 1)
$ws.Range("K1").Value = "use_as_trend"
$ws.Range("L1").Value = "match_year"

# New data cells (row 2)
$ws.Range("K2").Value = "T"
$ws.Range("L2").Value = 1953

# Column J (10th column) gets a best-fit custom width
$ws.Columns.Item(10).ColumnWidth = 19.25

# Make the sheet active and move the view/selection to match the target state
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J13").Select() | Out-Null
